$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their text formatting (values are stored as text, e.g. "302.93", "2.08%")
$ws.Range("D2:E27").NumberFormat = "@"
$ws.Range("D39:E50").NumberFormat = "@"

$ws.Range("D2").Value = "302.93"
$ws.Range("E2").Value = "2.08%"
$ws.Range("D3").Value = "32.09"
$ws.Range("E3").Value = "2.37%"
$ws.Range("D4").Value = "5.125"
$ws.Range("E4").Value = "0.60%"
$ws.Range("D5").Value = "0.07841"
$ws.Range("E5").Value = "-2.20%"
$ws.Range("D6").Value = "2.256"
$ws.Range("E6").Value = "-13.41%"
$ws.Range("D7").Value = "7.830"
$ws.Range("E7").Value = "0.83%"
$ws.Range("E8").Value = "0.56%"
$ws.Range("D9").Value = "0.9256"
$ws.Range("E9").Value = "-0.18%"
$ws.Range("E10").Value = "0.97%"
$ws.Range("D11").Value = "0.07811"
$ws.Range("E11").Value = "7.51%"
$ws.Range("D12").Value = "0.08877"
$ws.Range("E12").Value = "-0.50%"
$ws.Range("D13").Value = "0.03099"
$ws.Range("E13").Value = "2.69%"
$ws.Range("E14").Value = "0.33%"
$ws.Range("D15").Value = "0.001511"
$ws.Range("E15").Value = "1.10%"
$ws.Range("D16").Value = "0.005850"
$ws.Range("E16").Value = "-0.40%"
$ws.Range("D17").Value = "3.458"
$ws.Range("E17").Value = "-1.86%"
$ws.Range("D18").Value = "2.248"
$ws.Range("E18").Value = "0.06%"
$ws.Range("E19").Value = "1.28%"
$ws.Range("D20").Value = "0.1330"
$ws.Range("E20").Value = "-1.02%"
$ws.Range("D21").Value = "4.259"
$ws.Range("E21").Value = "25.89%"
$ws.Range("D22").Value = "0.1798"
$ws.Range("E22").Value = "8.92%"
$ws.Range("D23").Value = "0.04601"
$ws.Range("E23").Value = "0.33%"
$ws.Range("D24").Value = "0.001254"
$ws.Range("E24").Value = "1.05%"
$ws.Range("D25").Value = "0.004498"
$ws.Range("E25").Value = "1.73%"
$ws.Range("D26").Value = "0.0001253"
$ws.Range("E26").Value = "4.54%"
$ws.Range("E27").Value = "-1.06%"
$ws.Range("D39").Value = "0.01787"
$ws.Range("E39").Value = "1.66%"
$ws.Range("D40").Value = "0.04787"
$ws.Range("E40").Value = "7.31%"
$ws.Range("D41").Value = "0.007211"
$ws.Range("E41").Value = "5.09%"
$ws.Range("D42").Value = "0.1373"
$ws.Range("E42").Value = "2.50%"
$ws.Range("D43").Value = "0.002195"
$ws.Range("E43").Value = "2.50%"
$ws.Range("D44").Value = "0.009946"
$ws.Range("E44").Value = "3.96%"
$ws.Range("D45").Value = "0.00006259"
$ws.Range("E45").Value = "-4.42%"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("E46").Value = "0.36%"
$ws.Range("D47").Value = "0.003599"
$ws.Range("E47").Value = "-58.76%"
$ws.Range("D48").Value = "1.157"
$ws.Range("E48").Value = "40.99%"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("E49").Value = "0.36%"
$ws.Range("D50").Value = "0.0002002"
$ws.Range("E50").Value = "0.36%"
